# Atualização de bases das ligas, do dia: 06-04-2024 às 15:39
#
# This script applies an odds/fixture-data refresh to the
# "Uruguay Primera División" sheet:
#   - Rows 114/115 swap their match data (id + odds columns B,F:AC),
#     keeping the row-local A/C/D/E (row index, Div, Div Original Name,
#     Date) untouched.
#   - Rows 117/118/119 rotate their match data the same way
#     (117<-119, 118<-117, 119<-118).
#   - Rows 169-172 get refreshed odds; 169/170 keep the same fixture
#     (id/date/teams/opening odds) but with updated live odds, while
#     171/172 become brand-new upcoming fixtures. Row 173 only gets a
#     small tweak to its handicap odds (R/S).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 114 <-> Row 115 swap -------------------------------------------
$ws.Range("B114").Value = 7559469
$ws.Range("F114").Value = 'Montevideo Wanderers'
$ws.Range("G114").Value = 'Penarol'
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 'D'
$ws.Range("K114").Value = 4.75
$ws.Range("L114").Value = 3.4
$ws.Range("M114").Value = 1.7
$ws.Range("N114").Value = 2.7
$ws.Range("O114").Value = 3.2
$ws.Range("P114").Value = 2.45
$ws.Range("Q114").Value = 0
$ws.Range("R114").Value = 2.05
$ws.Range("S114").Value = 1.8
$ws.Range("T114").Value = 2.5
$ws.Range("U114").Value = 1.975
$ws.Range("V114").Value = 1.875
$ws.Range("W114").Value = -1
$ws.Range("X114").Value = 2.2
$ws.Range("Y114").Value = -1
$ws.Range("Z114").Value = 0
$ws.Range("AA114").Value = -0
$ws.Range("AB114").Value = -1
$ws.Range("AC114").Value = 0.875

$ws.Range("B115").Value = 7559468
$ws.Range("F115").Value = 'Liverpool Montevideo'
$ws.Range("G115").Value = 'CA River Plate'
$ws.Range("H115").Value = 2
$ws.Range("I115").Value = 1
$ws.Range("J115").Value = 'H'
$ws.Range("K115").Value = 1.7
$ws.Range("L115").Value = 3
$ws.Range("M115").Value = 5.75
$ws.Range("N115").Value = 1.833
$ws.Range("O115").Value = 3.2
$ws.Range("P115").Value = 4.5
$ws.Range("Q115").Value = -0.5
$ws.Range("R115").Value = 1.925
$ws.Range("S115").Value = 1.925
$ws.Range("T115").Value = 2.25
$ws.Range("U115").Value = 2.025
$ws.Range("V115").Value = 1.825
$ws.Range("W115").Value = 0.833
$ws.Range("X115").Value = -1
$ws.Range("Y115").Value = -1
$ws.Range("Z115").Value = 0.925
$ws.Range("AA115").Value = -1
$ws.Range("AB115").Value = 1.025
$ws.Range("AC115").Value = -1

# --- Row 117 <- 119, 118 <- 117, 119 <- 118 rotation ---------------------
$ws.Range("B117").Value = 7013409
$ws.Range("F117").Value = 'Nacional De Football'
$ws.Range("G117").Value = 'Torque'
$ws.Range("H117").Value = 1
$ws.Range("I117").Value = 1
$ws.Range("J117").Value = 'D'
$ws.Range("K117").Value = 1.666
$ws.Range("L117").Value = 3.9
$ws.Range("M117").Value = 4.5
$ws.Range("N117").Value = 1.615
$ws.Range("O117").Value = 4
$ws.Range("P117").Value = 4.75
$ws.Range("Q117").Value = -0.75
$ws.Range("R117").Value = 1.8
$ws.Range("S117").Value = 2.05
$ws.Range("T117").Value = 2.75
$ws.Range("U117").Value = 1.95
$ws.Range("V117").Value = 1.9
$ws.Range("W117").Value = -1
$ws.Range("X117").Value = 3
$ws.Range("Y117").Value = -1
$ws.Range("Z117").Value = -1
$ws.Range("AA117").Value = 1.05
$ws.Range("AB117").Value = -1
$ws.Range("AC117").Value = 0.8999999999999999

$ws.Range("B118").Value = 7013886
$ws.Range("F118").Value = 'Racing Club de Montevideo'
$ws.Range("G118").Value = 'Cerro'
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 1
$ws.Range("J118").Value = 'A'
$ws.Range("K118").Value = 2.25
$ws.Range("L118").Value = 3.1
$ws.Range("M118").Value = 3.25
$ws.Range("N118").Value = 2.25
$ws.Range("O118").Value = 2.875
$ws.Range("P118").Value = 3.5
$ws.Range("Q118").Value = -0.25
$ws.Range("R118").Value = 1.95
$ws.Range("S118").Value = 1.9
$ws.Range("T118").Value = 2
$ws.Range("U118").Value = 1.925
$ws.Range("V118").Value = 1.925
$ws.Range("Y118").Value = 2.5
$ws.Range("Z118").Value = -1
$ws.Range("AA118").Value = 0.8999999999999999
$ws.Range("AB118").Value = -1
$ws.Range("AC118").Value = 0.925

$ws.Range("B119").Value = 7013885
$ws.Range("F119").Value = 'La Luz'
$ws.Range("G119").Value = 'Atletico Fenix Montevideo'
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 2
$ws.Range("J119").Value = 'A'
$ws.Range("K119").Value = 3
$ws.Range("L119").Value = 3
$ws.Range("M119").Value = 2.4
$ws.Range("N119").Value = 2.9
$ws.Range("O119").Value = 2.75
$ws.Range("P119").Value = 2.6
$ws.Range("Q119").Value = 0
$ws.Range("R119").Value = 2.025
$ws.Range("S119").Value = 1.825
$ws.Range("T119").Value = 2
$ws.Range("U119").Value = 2.025
$ws.Range("V119").Value = 1.825
$ws.Range("W119").Value = -1
$ws.Range("X119").Value = -1
$ws.Range("Y119").Value = 1.6
$ws.Range("Z119").Value = -1
$ws.Range("AA119").Value = 0.825
$ws.Range("AB119").Value = 0
$ws.Range("AC119").Value = -0

# --- Rows 169-172 refreshed fixtures / odds, row 173 odds tweak ----------
$ws.Range("B169").Value = 8014043
$ws.Range("E169").Value = 45388.625
$ws.Range("F169").Value = 'Danubio'
$ws.Range("G169").Value = 'Club Atletico Progreso'
$ws.Range("K169").Value = 2.375
$ws.Range("L169").Value = 3.1
$ws.Range("M169").Value = 3
$ws.Range("N169").Value = 2.25
$ws.Range("O169").Value = 3.1
$ws.Range("P169").Value = 3.25
$ws.Range("Q169").Value = -0.25
$ws.Range("R169").Value = 1.95
$ws.Range("S169").Value = 1.9
$ws.Range("T169").Value = 2.25
$ws.Range("U169").Value = 1.925
$ws.Range("V169").Value = 1.925

$ws.Range("B170").Value = 8014089
$ws.Range("E170").Value = 45388.75
$ws.Range("F170").Value = 'Nacional De Football'
$ws.Range("G170").Value = 'Cerro Largo'
$ws.Range("K170").Value = 1.5
$ws.Range("L170").Value = 4
$ws.Range("M170").Value = 6
$ws.Range("N170").Value = 1.6
$ws.Range("O170").Value = 3.8
$ws.Range("P170").Value = 5
$ws.Range("Q170").Value = -0.75
$ws.Range("R170").Value = 1.825
$ws.Range("S170").Value = 2.025
$ws.Range("T170").Value = 2.25
$ws.Range("U170").Value = 1.85
$ws.Range("V170").Value = 2

$ws.Range("B171").Value = 8014090
$ws.Range("E171").Value = 45389.375
$ws.Range("F171").Value = 'Atletico Fenix Montevideo'
$ws.Range("G171").Value = 'Liverpool Montevideo'
$ws.Range("K171").Value = 2.625
$ws.Range("L171").Value = 2.9
$ws.Range("M171").Value = 2.9
$ws.Range("N171").Value = 2.625
$ws.Range("O171").Value = 2.9
$ws.Range("P171").Value = 2.875
$ws.Range("Q171").Value = 0
$ws.Range("R171").Value = 1.85
$ws.Range("S171").Value = 2
$ws.Range("T171").Value = 2
$ws.Range("U171").Value = 1.875
$ws.Range("V171").Value = 1.975

$ws.Range("B172").Value = 8014132
$ws.Range("E172").Value = 45389.41666666666
$ws.Range("F172").Value = 'Cerro'
$ws.Range("G172").Value = 'Rampla Juniors'
$ws.Range("K172").Value = 2.2
$ws.Range("L172").Value = 3.2
$ws.Range("M172").Value = 3.2
$ws.Range("N172").Value = 2.2
$ws.Range("O172").Value = 3.2
$ws.Range("P172").Value = 3.2
$ws.Range("Q172").Value = -0.25
$ws.Range("R172").Value = 1.95
$ws.Range("S172").Value = 1.9
$ws.Range("T172").Value = 2.25
$ws.Range("U172").Value = 1.95
$ws.Range("V172").Value = 1.9

$ws.Range("R173").Value = 1.875
$ws.Range("S173").Value = 1.975
